$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text so numeric-looking strings
# (e.g. "1.002") are not auto-converted to numbers on assignment,
# matching the source data which stores these as plain strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.180.72'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.906.26'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '326.06'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.3892'
$ws.Range("E8").Value = '  -1.33%  '
$ws.Range("D9").Value = '0.07886'
$ws.Range("E9").Value = '  -0.76%  '
$ws.Range("D10").Value = '0.9916'
$ws.Range("E10").Value = '  -1.23%  '
$ws.Range("D11").Value = '21.96'
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").Value = '1.898.14'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").Value = '5.775'
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("D14").Value = '7.054'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '88.12'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '0.000009945'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("D19").Value = '17.09'
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("D21").Value = '29.169.39'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '5.317'
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("D23").Value = '11.16'
$ws.Range("E23").Value = '  +0.40%  '
$ws.Range("D24").Value = '2.105'
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").Value = '156.21'
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").Value = '19.43'
$ws.Range("E26").Value = '  -0.37%  '
$ws.Range("D27").Value = '5.926'
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("D28").Value = '118.71'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '1.890'
$ws.Range("E29").Value = '  -5.47%  '
$ws.Range("D30").Value = '0.09353'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '0.8966'
$ws.Range("E31").Value = '  -3.46%  '
$ws.Range("D32").Value = '5.248'
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("E33").Value = '  -2.49%  '
$ws.Range("D34").Value = '3.164'
$ws.Range("E34").Value = '  -3.31%  '
$ws.Range("D35").Value = '0.05804'
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").Value = '1.175'
$ws.Range("E36").Value = '  -2.70%  '
$ws.Range("D37").Value = '0.02088'
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '1.001'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '7.682'
$ws.Range("E39").Value = '  -3.15%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5701'
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("D41").Value = '0.1799'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = '9.732'
$ws.Range("E42").Value = '  -2.50%  '
$ws.Range("D43").Value = '11.92'
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").Value = '0.5353'
$ws.Range("E44").Value = '  -1.11%  '
$ws.Range("D45").Value = '2.183'
$ws.Range("E45").Value = '  -3.38%  '
$ws.Range("D46").Value = '0.07018'
$ws.Range("D47").Value = '1.851'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = '2.551'
$ws.Range("E48").Value = '  -0.14%  '
$ws.Range("D49").Value = '113.19'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("D50").Value = '0.2939'
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '1.037'
$ws.Range("E51").Value = '  -2.73%  '

# Restore column D formatting/style to its original (General) state.
$ws.Range("D2:D51").Style = "Normal"

